# Organização do código fonte.
# - Corrige o texto "ROBERTa" para "RoBERTa" na célula A6
# - Adiciona uma nova linha (A7) com a mesma formatação (negrito) das demais
#   células da coluna A, preparando a planilha para um novo modelo
# - Move a seleção ativa para a nova célula A7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Corrige a capitalização de "ROBERTa" para "RoBERTa"
$ws.Range("A6").Value = "RoBERTa"

# Usa a célula A6 (negrito) como modelo de estilo para a nova célula A7
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A7").Value = $null

$excel.CutCopyMode = 0

# Atualiza a seleção ativa para a nova célula
$ws.Range("A7").Select() | Out-Null
